$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.866.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "'3.529.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.15%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'613.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").Value = "'173.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("E7").Value = "  -1.39%  "

$ws.Range("D8").Value = "'3.522.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.17%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("D11").Value = "'7.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "

$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").Value = "'4.099.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.09%  "

$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").Value = "'614.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").Value = "'3.531.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").Value = "'70.857.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("D21").Value = "'17.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.85%  "

$ws.Range("D22").Value = "'0.886"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "'9.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.09%  "

$ws.Range("D24").Value = "'15.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "

$ws.Range("D25").Value = "'98.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("E26").Value = "  -1.60%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  -0.65%  "

$ws.Range("D29").Value = "'33.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.56%  "

$ws.Range("D30").Value = "'9.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("D32").Value = "'8.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.11%  "

$ws.Range("E33").Value = "  -0.51%  "

$ws.Range("E34").Value = "  -1.41%  "

$ws.Range("D35").Value = "'605.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.45%  "

$ws.Range("E36").Value = "  -0.71%  "

$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").Value = "'3.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.55%  "

$ws.Range("D39").Value = "'0.0471"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("D40").Value = "'57.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("E43").Value = "  +5.44%  "

$ws.Range("D44").Value = "'3.365.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("E45").Value = "  -2.18%  "

$ws.Range("D46").Value = "'2.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").Value = "'32.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.45%  "

$ws.Range("E48").Value = "  -1.88%  "

$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").Value = "'133.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.22%  "
